$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths / hidden helper columns (cols N..V = 14..22)
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 10.0                # N -> width ~10.77734375
$ws.Columns.Item(15).ColumnWidth = 16.333333333333332  # O -> width ~17.21875
$ws.Columns.Item(16).ColumnWidth = 26.166666666666668  # P -> width 27
$ws.Columns.Item(17).ColumnWidth = 8.0                 # Q -> width ~8.88671875
$ws.Columns.Item(18).ColumnWidth = 8.0                 # R -> width ~8.88671875
$ws.Columns.Item(19).ColumnWidth = 1.3333333333333333  # S -> width ~2.109375
$ws.Columns.Item(20).ColumnWidth = 18.833333333333332  # T -> width ~19.6640625
$ws.Columns.Item(21).ColumnWidth = 20.333333333333332  # U -> width ~21.109375
$ws.Columns.Item(22).ColumnWidth = 22.333333333333332  # V -> width ~23.109375

$ws.Columns.Item(17).Hidden = $true  # Q
$ws.Columns.Item(18).Hidden = $true  # R
$ws.Columns.Item(19).Hidden = $true  # S

# ---------------------------------------------------------------------------
# Row 1 header values (right-to-left sheet, author filled from V back to G)
# ---------------------------------------------------------------------------
$ws.Range("V1").Value = "release version"
$ws.Range("U1").Value = "task name"
$ws.Range("T1").Value = "task assigned to "
$ws.Range("P1").Value = "Task Reviewed By"

# Style the P1/T1/U1/V1 headers: bold, 14pt, centered
foreach ($addr in @("P1", "T1", "U1", "V1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Font.Size = 14
    $c.HorizontalAlignment = -4108
}

# Merge G1:O1 and style the big "Comments" title (bold, 16pt, centered)
$title = $ws.Range("G1:O1")
$title.Font.Bold = $true
$title.Font.Size = 16
$title.HorizontalAlignment = -4108
$title.Merge()
$ws.Range("G1").Value = "Comments"
# The merged cell anchor ends up slightly smaller (14pt) than the rest of the
# merge range (16pt), matching the author's final formatting pass.
$ws.Range("G1").Font.Size = 14

$ws.Rows.Item(1).RowHeight = 21

# ---------------------------------------------------------------------------
# Row 2 values
# ---------------------------------------------------------------------------
$ws.Range("V2").Value = "V1"
$ws.Range("U2").Value = "SRS_adminFeatures"
$ws.Range("T2").Value = " Dina"
$ws.Range("P2").Value = "Mayar"
$ws.Range("O2").Value = "clarify on (SRS_admin_editCst_02,SRS_admin_editAcc_05) what is the fields the admain which can the admin edit "

foreach ($addr in @("P2", "T2", "U2", "V2")) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# View state
# ---------------------------------------------------------------------------
$ws.Range("P16").Select()

Write-Output "done"
